$p = $ppt.ActivePresentation

# Duplicate slide 6 (DiSCoVER: top drugs table) and move the duplicate to the end (slide 9)
$src = $p.Slides.Item(6)
$newRange = $src.Duplicate()
$new = $newRange.Item(1)
$new.MoveTo($p.Slides.Count)
